# teste conversao de texto pra numero
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("produtos")

$data = @(
    @("cookies limão 140g", "R$20.99", "R$19.99"),
    @("cookies 3 chocolates 140g", "R$20.99", "R$19.99"),
    @("cookies com gotas ao leite 120g", "R$17.99", "R$16.99"),
    @("wafer mil folhas 80g", "R$17.99", "R$15.99"),
    @("wafer lanut 80g", "R$17.99", "R$15.99"),
    @("wafer cappuccino 80g", "R$17.99", "R$15.99"),
    @("biscoito caramelizado coberto 100g", "R$17.99", "R$15.99"),
    @("cookies duo chocolate 120g", "R$17.99", "R$16.99"),
    @("biscoito tradicional ao leite 100g", "R$17.99", "R$15.99"),
    @("wafer limão 80g", "R$17.99", "R$15.99"),
    @("wafer coberto ao leite 80g", "R$17.99", "R$16.99"),
    @("pão de mel doce de leite 50g", "R$11.99", "R$9.99"),
    @("biscoito tradicional coberto ao leite 50g", "R$10.99", "R$9.99"),
    @("cookies tradicional gotas ao leite 40g", "R$10.99", "R$9.99"),
    @("biscoito caramelizado 100g", "R$10.99", "R$9.99"),
    @("biscoito caramelizado coberto ao leite 50g", "R$10.99", "R$9.99"),
    @("bytes wafer avelã 50g", "R$10.59", "R$9.99"),
    @("wafer bytes sintonia 45g", "R$10.59", "R$9.99")
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}
